# B6-PowerPoint.pptx edit
#
# 1. Re-style the three tables (slides 14, 15, 16) from the local
#    "Table_0" style {66136E33-A210-4CF8-AA54-9CA5D72045B3} to the
#    built-in table style {1079D921-6634-4932-8769-F80908AC00E2}
#    (picked from the Table Tools > Design > Table Styles gallery).
#
# 2. Swap the presentation's applied theme color scheme from the
#    "Integral" / "Red Violet" palette back to the default
#    "Office" palette (Design > Themes > Office) -- i.e. the colours
#    that live in the theme part feeding the (only) slide master.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Part 1: table styles
# ---------------------------------------------------------------
$targetStyleId = "{1079D921-6634-4932-8769-F80908AC00E2}"

foreach ($slideNum in 14, 15, 16) {
    $slide = $p.Slides.Item($slideNum)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($targetStyleId)
        }
    }
}

# ---------------------------------------------------------------
# Part 2: theme colors -> back to stock "Office" scheme
# ---------------------------------------------------------------
function HexToOle([string]$hex) {
    $rr = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $gg = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $bb = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $bb * 65536 + $gg * 256 + $rr
}

# Order matches ThemeColorScheme.Colors(1..12):
#  1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#  8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $colorScheme.Colors($i).RGB = HexToOle $officeColors[$i - 1]
}
